$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.219.36"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.648.10"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.878.99"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "1.645.03"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "27.210.86"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").Value = "1.265.42"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  +4.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "1.788.99"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.71%  "
